$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.943.88'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '2.486.88'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.516'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.532'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.98%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('E12').Value = '  -2.90%  '
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('D15').Value = '2.874.78'
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('D16').Value = '2.490.47'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.823'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.41%  '
$ws.Range('D18').Value = '47.736.07'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('E19').Value = '  +8.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.55'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').Value = '0.0₃0924'
$ws.Range('E22').Value = '  -2.45%  '
$ws.Range('E23').Value = '  -2.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.99%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('E28').Value = '  +4.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.64'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.15%  '
$ws.Range('E30').Value = '  -6.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.26'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.88'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.61%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.26'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0768'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.54'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.35%  '
$ws.Range('E39').Value = '  -5.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '121.92'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.40%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.110'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.20'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.35%  '
$ws.Range('E43').Value = '  +0.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0299'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = '1.990.69'
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('E46').Value = '  -2.27%  '
$ws.Range('E47').Value = '  -1.08%  '
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.38'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.40%  '
